$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.051.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.67%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.621.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.44%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.518"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.75%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  -1.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0625"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.82%  "

$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.626.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.05%  "

$ws.Range("E13").Value = "  -0.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.542"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.020.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0742"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "216.10"
$ws.Range("D18").Style = "Normal"

$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.53%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.56%  "

$ws.Range("E27").Value = "  -0.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0505"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("E30").Value = "  -1.20%  "

$ws.Range("E31").Value = "  -1.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.340.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.29%  "

$ws.Range("E34").Value = "  -0.50%  "

$ws.Range("E35").Value = "  -0.38%  "

$ws.Range("E36").Value = "  -1.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.546"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.851"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.22%  "

$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.802"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.83%  "

$ws.Range("E41").Value = "  -0.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.757.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.48%  "

$ws.Range("E44").Value = "  -2.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.35%  "

$ws.Range("E46").Value = "  +0.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.855"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +28.81%  "

$ws.Range("E48").Value = "  -1.01%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0512"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.39%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0995"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.60%  "

Write-Host "Applied cryptos list update"